# Update "想去人数" (want-to-go count) values in column F across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 18
$ws.Range("F3").Value = 984
$ws.Range("F4").Value = 225
$ws.Range("F6").Value = 1129
$ws.Range("F7").Value = 906
$ws.Range("F8").Value = 280
$ws.Range("F11").Value = 881
$ws.Range("F12").Value = 310
$ws.Range("F13").Value = 593
$ws.Range("F14").Value = 519
$ws.Range("F17").Value = 1249
$ws.Range("F18").Value = 2928
$ws.Range("F20").Value = 1526
$ws.Range("F21").Value = 1291
$ws.Range("F26").Value = 1058
$ws.Range("F28").Value = 3258
$ws.Range("F29").Value = 636
$ws.Range("F31").Value = 1451

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 17
$ws.Range("F4").Value = 14
$ws.Range("F5").Value = 63
$ws.Range("F7").Value = 165
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = 35
$ws.Range("F10").Value = 18
$ws.Range("F12").Value = 25

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 768

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 18
$ws.Range("F3").Value = 768
$ws.Range("F4").Value = 17
$ws.Range("F6").Value = 984
$ws.Range("F7").Value = 225
$ws.Range("F8").Value = 14
$ws.Range("F10").Value = 1129
$ws.Range("F11").Value = 906
$ws.Range("F12").Value = 280
$ws.Range("F13").Value = 63
$ws.Range("F17").Value = 165
$ws.Range("F18").Value = 6
$ws.Range("F19").Value = 35
$ws.Range("F21").Value = 18
$ws.Range("F22").Value = 881
$ws.Range("F23").Value = 310
$ws.Range("F24").Value = 593
$ws.Range("F25").Value = 519
$ws.Range("F28").Value = 1249
$ws.Range("F29").Value = 2928
$ws.Range("F31").Value = 1526
$ws.Range("F32").Value = 1291
$ws.Range("F38").Value = 25
$ws.Range("F39").Value = 1058
$ws.Range("F41").Value = 3258
$ws.Range("F42").Value = 636
$ws.Range("F44").Value = 1451
